# Refine CPS (Cognitive Spontaneity) barrier/question wording so it maps
# precisely to items 9 (role variety) and 10 (unconventional object use),
# replacing the vague "repeating storylines" framing.

$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2 (used as literal values below, as
# the PowerShell COM shim exposes Find.Execute positionally like VBA does).

# 1) "Barrier" cell text (table summarising facilitator observations).
$d.Content.Find.Execute(
    "孩子等待指令才行動；每日重複同一故事線，缺乏新意", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "孩子等待指令才行動，不主動創作玩法；孩子只扮演固定角色，鮮少以創意方式運用物料", `
    2) | Out-Null

# 2) "Reflection question" cell text directly following the barrier cell.
$d.Content.Find.Execute(
    "如何延伸孩子的想法而不變成教師主導的教學活動？", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "如何支持孩子自主創作玩法、嘗試不同角色，並以新方式運用物料，而非依賴教師指令或固定劇情？", `
    2) | Out-Null

# 3) Middle line of the multi-line facilitator quote (line separated from
#    its neighbours by literal newline characters inside the same run).
$d.Content.Find.Execute(
    "每日都玩同一條故事線，冇乜新變化。", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "每日扮演相同角色（永遠係廚師／媽媽），唔嘗試新角色；或只按物品既有用法玩，唔作創意延伸。", `
    2) | Out-Null

Write-Output "Done."
